# Correct status name labels in the "statut_label" (column B) and
# "statut_name" (column C) columns.
#
# - "bleu" -> "noir"
# - "pas de résultat ni de publication" -> "pas de résultat postés ni publiés"
# - "résultat et / ou publication posté" -> "résultat postés ou publiés"
# - "résultat et / ou publication posté dans les 36 mois" -> "résultat postés ou publiés dans les 36 mois"
# - "résultat et / ou publication posté dans les 12 mois" -> "résultat postés ou publiés dans les 12 mois"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Longer / more specific strings first so that a shorter match (e.g. the
# "pas de résultat et / ou publication posté" prefix) doesn't get replaced
# before the longer variants are handled.
$ws.Cells.Replace("résultat et / ou publication posté dans les 36 mois", "résultat postés ou publiés dans les 36 mois")
$ws.Cells.Replace("résultat et / ou publication posté dans les 12 mois", "résultat postés ou publiés dans les 12 mois")
$ws.Cells.Replace("résultat et / ou publication posté", "résultat postés ou publiés")
$ws.Cells.Replace("pas de résultat ni de publication", "pas de résultat postés ni publiés")
$ws.Cells.Replace("bleu", "noir")
